$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 20, shifting existing rows 20-26 down to 21-27
$ws.Rows.Item(20).Insert()

# Populate the new row 20 with the new weekly price entry
$ws.Cells.Item(20, 1).Value2 = 5
$ws.Cells.Item(20, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(20, 3).Value = "Maule"
$ws.Cells.Item(20, 4).Value2 = 44518
$ws.Cells.Item(20, 5).Value2 = 7
$ws.Cells.Item(20, 6).Value = "Fruta"
$ws.Cells.Item(20, 7).Value2 = 100101
$ws.Cells.Item(20, 8).Value = "Berries"
$ws.Cells.Item(20, 9).Value2 = 100101001
$ws.Cells.Item(20, 10).Value = "Arándano (blue)"
$ws.Cells.Item(20, 11).Value = "Sin especificar"
$ws.Cells.Item(20, 12).Value = "Primera"
$ws.Cells.Item(20, 13).Value2 = 20
$ws.Cells.Item(20, 14).Value2 = 5000
$ws.Cells.Item(20, 15).Value2 = 5000
$ws.Cells.Item(20, 16).Value2 = 5000
$ws.Cells.Item(20, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(20, 18).Value = "Provincia de Linares"
$ws.Cells.Item(20, 19).Value2 = 2500
$ws.Cells.Item(20, 20).Value2 = 2

$ws.Cells.Item(20, 4).NumberFormat = $ws.Cells.Item(21, 4).NumberFormat
